$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2522.75
$ws.Range("I40").Value = 2400.0625
$ws.Range("J40").Value = 2768.125
$ws.Range("K40").Value = 2400.0625
$ws.Range("L40").Value = 2768.125
$ws.Range("M40").Value = -2225.0625
$ws.Range("N40").Value = -3118.125

$ws.Range("H80").Value = 44124.695
$ws.Range("I80").Value = 374.8889
$ws.Range("J80").Value = 72249.57000000001
$ws.Range("K80").Value = 1124.6667
$ws.Range("L80").Value = 216748.71
$ws.Range("M80").Value = -126.6667
$ws.Range("N80").Value = -218744.71

$ws.Range("H83").Value = 44124.695
$ws.Range("I83").Value = 374.8889
$ws.Range("J83").Value = 72249.57000000001
$ws.Range("K83").Value = 3374.0001
$ws.Range("L83").Value = 650246.1300000001
$ws.Range("M83").Value = 1617.9999
$ws.Range("N83").Value = -660230.1300000001

$ws.Range("H86").Value = 10529827
$ws.Range("I86").Value = 1889.4445
$ws.Range("J86").Value = 20004970
$ws.Range("K86").Value = 1889.4445
$ws.Range("L86").Value = 20004970
$ws.Range("M86").Value = -766.4445000000001
$ws.Range("N86").Value = -20007216

$ws.Range("H89").Value = 10529827
$ws.Range("I89").Value = 1889.4445
$ws.Range("J89").Value = 20004970
$ws.Range("K89").Value = 9447.2225
$ws.Range("L89").Value = 100024850
$ws.Range("M89").Value = -3831.2225
$ws.Range("N89").Value = -100036082

$ws.Range("H101").Value = 1456.9
$ws.Range("I101").Value = 336.66666
$ws.Range("J101").Value = 3137.25
$ws.Range("K101").Value = 1009.99998
$ws.Range("L101").Value = 9411.75
$ws.Range("M101").Value = 612.0000200000001
$ws.Range("N101").Value = -12655.75

$ws.Range("H128").Value = 27964.334
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 27964.334
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 27964.334
$ws.Range("M128").Value = ""
$ws.Range("N128").Value = -37924.334

$ws.Range("H136").Value = 54930
$ws.Range("J136").Value = 54930
$ws.Range("L136").Value = 54930
$ws.Range("N136").Value = -65130

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 39800
$ws.Range("J15").Value = 39800
$ws.Range("L15").Value = 39800
$ws.Range("N15").Value = -40500

$ws.Range("H32").Value = 22202.215
$ws.Range("I32").Value = 4779.7534
$ws.Range("K32").Value = 4779.7534
$ws.Range("M32").Value = -4492.7534

$ws.Range("H97").Value = 44439.043
$ws.Range("I97").Value = 46322.637
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 46322.637
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -45826.637
$ws.Range("N97").Value = -3992

$ws.Range("H134").Value = 67424.5
$ws.Range("J134").Value = 67424.5
$ws.Range("L134").Value = 67424.5
$ws.Range("N134").Value = -77564.5

$ws.Range("H135").Value = 34965.75
$ws.Range("J135").Value = 34965.75
$ws.Range("L135").Value = 34965.75
$ws.Range("N135").Value = -45105.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 398000
$ws.Range("J43").Value = 398000
$ws.Range("L43").Value = 398000
$ws.Range("N43").Value = -398362

$ws.Range("H63").Value = 45265.5
$ws.Range("J63").Value = 45265.5
$ws.Range("L63").Value = 45265.5
$ws.Range("N63").Value = -46637.5

$ws.Range("H66").Value = 45265.5
$ws.Range("J66").Value = 45265.5
$ws.Range("L66").Value = 135796.5
$ws.Range("N66").Value = -142660.5

$ws.Range("H94").Value = 824
$ws.Range("I94").Value = 707.5
$ws.Range("J94").Value = 901.6667
$ws.Range("K94").Value = 707.5
$ws.Range("L94").Value = 901.6667
$ws.Range("M94").Value = -256.5
$ws.Range("N94").Value = -1803.6667

$ws.Range("H137").Value = 39983
$ws.Range("J137").Value = 39983
$ws.Range("L137").Value = 39983
$ws.Range("N137").Value = -50183

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28079.686
$ws.Range("I31").Value = 1286.7916
$ws.Range("J31").Value = 49514
$ws.Range("K31").Value = 1286.7916
$ws.Range("L31").Value = 49514
$ws.Range("M31").Value = -991.7916
$ws.Range("N31").Value = -50104

$ws.Range("H34").Value = 28079.686
$ws.Range("I34").Value = 1286.7916
$ws.Range("J34").Value = 49514
$ws.Range("K34").Value = 1286.7916
$ws.Range("L34").Value = 49514
$ws.Range("M34").Value = -1084.7916
$ws.Range("N34").Value = -49918

$ws.Range("H132").Value = 26318572
$ws.Range("I132").Value = 23258456
$ws.Range("J132").Value = 35717500
$ws.Range("K132").Value = 69775368
$ws.Range("L132").Value = 107152500
$ws.Range("M132").Value = -69772838
$ws.Range("N132").Value = -107157560

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 746
$ws.Range("I2").Value = 199.33333
$ws.Range("J2").Value = 1074
$ws.Range("K2").Value = 1195.99998
$ws.Range("L2").Value = 6444
$ws.Range("M2").Value = -1082.99998
$ws.Range("N2").Value = -6670

$ws.Range("H58").Value = 2839.8
$ws.Range("I58").Value = 1200
$ws.Range("J58").Value = 3249.75
$ws.Range("K58").Value = 3600
$ws.Range("L58").Value = 9749.25
$ws.Range("M58").Value = -3472
$ws.Range("N58").Value = -10005.25

$ws.Range("H122").Value = 367.33334
$ws.Range("I122").Value = 367.33334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 3306.00006
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -856.0000600000003
$ws.Range("N122").Value = ""

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 77003900
$ws.Range("I80").Value = 111226950
$ws.Range("J80").Value = 2050
$ws.Range("K80").Value = 111226950
$ws.Range("L80").Value = 2050
$ws.Range("M80").Value = -111225952
$ws.Range("N80").Value = -4046

$ws.Range("H83").Value = 77003900
$ws.Range("I83").Value = 111226950
$ws.Range("J83").Value = 2050
$ws.Range("K83").Value = 556134750
$ws.Range("L83").Value = 10250
$ws.Range("M83").Value = -556129758
$ws.Range("N83").Value = -20234

$ws.Range("H112").Value = 39999
$ws.Range("J112").Value = 39999
$ws.Range("L112").Value = 39999
$ws.Range("N112").Value = -42215

$ws.Range("H132").Value = 3803.7856
$ws.Range("I132").Value = 3416.25
$ws.Range("K132").Value = 10248.75
$ws.Range("M132").Value = -7718.75

$ws.Range("H134").Value = 30333.133
$ws.Range("J134").Value = 30333.133
$ws.Range("L134").Value = 90999.399
$ws.Range("N134").Value = -96069.399

$ws.Range("H136").Value = 20458.555
$ws.Range("J136").Value = 20458.555
$ws.Range("L136").Value = 61375.665
$ws.Range("N136").Value = -66475.66500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 39344.31
$ws.Range("I16").Value = 48141.477
$ws.Range("J16").Value = 2396.2
$ws.Range("K16").Value = 48141.477
$ws.Range("L16").Value = 2396.2
$ws.Range("M16").Value = -47971.477
$ws.Range("N16").Value = -2736.2

$ws.Range("H55").Value = 393.34482
$ws.Range("I55").Value = 240.84616
$ws.Range("J55").Value = 517.25
$ws.Range("K55").Value = 240.84616
$ws.Range("L55").Value = 517.25
$ws.Range("M55").Value = -67.84616
$ws.Range("N55").Value = -863.25

$ws.Range("H88").Value = 18975
$ws.Range("I88").Value = 6917.75
$ws.Range("J88").Value = 43089.5
$ws.Range("K88").Value = 6917.75
$ws.Range("L88").Value = 43089.5
$ws.Range("M88").Value = -6489.75
$ws.Range("N88").Value = -43945.5

$ws.Range("H91").Value = 18975
$ws.Range("I91").Value = 6917.75
$ws.Range("J91").Value = 43089.5
$ws.Range("K91").Value = 6917.75
$ws.Range("L91").Value = 43089.5
$ws.Range("M91").Value = -5435.75
$ws.Range("N91").Value = -46053.5

$ws.Range("H92").Value = 24999.334
$ws.Range("J92").Value = 24999.334
$ws.Range("L92").Value = 24999.334
$ws.Range("N92").Value = -29991.334

$ws.Range("H93").Value = 1661.7273
$ws.Range("J93").Value = 949.5
$ws.Range("L93").Value = 949.5
$ws.Range("N93").Value = -3445.5

$ws.Range("H132").Value = 3786.0386
$ws.Range("I132").Value = 3660.4348
$ws.Range("J132").Value = 4749
$ws.Range("K132").Value = 10981.3044
$ws.Range("L132").Value = 14247
$ws.Range("M132").Value = -8451.304400000001
$ws.Range("N132").Value = -19307

$ws.Range("H134").Value = 64925.5
$ws.Range("J134").Value = 64925.5
$ws.Range("L134").Value = 64925.5
$ws.Range("N134").Value = -75065.5

$ws.Range("H135").Value = 35900
$ws.Range("J135").Value = 35900
$ws.Range("L135").Value = 35900
$ws.Range("N135").Value = -46040

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2121.918
$ws.Range("I132").Value = 2219.5652
$ws.Range("J132").Value = 1822.4667
$ws.Range("K132").Value = 6658.6956
$ws.Range("L132").Value = 5467.4001
$ws.Range("M132").Value = -4128.6956
$ws.Range("N132").Value = -10527.4001

$ws.Range("H135").Value = 40210.383
$ws.Range("J135").Value = 40210.383
$ws.Range("L135").Value = 40210.383
$ws.Range("N135").Value = -50350.383

$ws.Range("H137").Value = 32642.834
$ws.Range("J137").Value = 32642.834
$ws.Range("L137").Value = 32642.834
$ws.Range("N137").Value = -42842.834
